# "units out of raw detection template"
#
# Fills in row 3 (the live "worked example" row used to preview the
# N3/O3 label-building formulas) with sample data, and drops the
# "units" (column D) component out of the N3 label-concatenation
# formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sample data for the worked example in row 3.
$ws.Range("B3").Value = "stool"
$ws.Range("C3").Value = "TAC"
$ws.Range("D3").Value = "Ct value"
$ws.Range("E3").Value = "Virus"
$ws.Range("F3").Value = "Dengue"
$ws.Range("G3").Value = "Dengue"

# N3's label formula previously appended the "units" field ($D3) just
# before the trailing ", by <assay type> result" suffix. Take that
# units component back out of the concatenation.
$n3Formula = '=TRIM(IF($H3="",$G3,"")' + "`n" `
    + '&IF($H3<>"",$H3,"")' + "`n" `
    + '&IF($I3<>""," "&$I3,"")' + "`n" `
    + '&IF(OR($I3="LT",$I3="ST",AND($D3="",$I3<>"")),"-pos","")' + "`n" `
    + '&IF($K3<>""," "&$K3,"")' + "`n" `
    + '&IF($J3<>""," "&$J3&"-pos","")' + "`n" `
    + '&IF($L3<>""," "&$L3&"-neg","")' + "`n" `
    + '&", by "&$C3&" result")'
$ws.Range("N3").Formula = $n3Formula

# Re-select the formula cell that was being worked on (matches the
# saved sheetView selection).
$ws.Range("N3").Select()
